$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.910.03"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.246.79"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "'114.13"
$ws.Range("E5").Value = "  -3.03%  "
$ws.Range("D6").Value = "'300.42"
$ws.Range("E6").Value = "  +12.15%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "'45.87"
$ws.Range("E10").Value = "  -5.12%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "'56.35"
$ws.Range("E12").Value = "  +2.78%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'15.38"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "'0.900"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "2.586.45"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "2.269.73"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "42.831.41"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").Value = "'7.72"
$ws.Range("E20").Value = "  +10.47%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  +27.33%  "
$ws.Range("D23").Value = "'73.46"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("D25").Value = "'232.83"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'40.11"
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "'175.78"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'21.29"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "'0.0902"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  +12.12%  "
$ws.Range("E37").Value = "  +5.37%  "
$ws.Range("D38").Value = "'0.128"
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "'0.240"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "'72.11"
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("E44").Value = "  -6.15%  "
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "'5.62"
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("D48").Value = "'1.38"
$ws.Range("E48").Value = "  +7.22%  "
$ws.Range("D49").Value = "'106.19"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("D50").Value = "'8.72"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "'0.0989"
$ws.Range("E51").Value = "  -1.67%  "

# Restore default (unstyled) formatting for cells that were forced to text,
# so the leading apostrophe does not leave a text-format/quote-prefix style behind.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
